$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename existing "Complete" status to "Completed" (affects F2:F4 which share this string)
$ws2.Range("F2").Value = "Completed"
$ws2.Range("F3").Value = "Completed"
$ws2.Range("F4").Value = "Completed"

# Fill in previously-empty status cells with new statuses
$ws2.Range("F5").Value = "In progress"
$ws2.Range("F6").Value = "In progress"
$ws2.Range("F7").Value = "In progress"
$ws2.Range("F8").Value = "Not started"
$ws2.Range("F9").Value = "Not started"
$ws2.Range("F10").Value = "Not started"

# Update the active selection on Sheet2 to F8
$ws2.Activate() | Out-Null
$ws2.Range("F8").Select() | Out-Null
